$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 78, shifting existing rows 78..204 down to 79..205
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new weekly data point
$ws.Range("A78").Value = 3
$ws.Range("B78").Value = "Femacal de La Calera"
$ws.Range("C78").Value = "Coquimbo"
$ws.Range("D78").Value = 44557
$ws.Range("E78").Value = 5
$ws.Range("F78").Value = 100112001
$ws.Range("G78").Value = "Berenjena"
$ws.Range("H78").Value = "Sin especificar"
$ws.Range("I78").Value = "Primera"
$ws.Range("J78").Value = 50
$ws.Range("K78").Value = 8000
$ws.Range("L78").Value = 8000
$ws.Range("M78").Value = 8000
$ws.Range("N78").Value = "`$/caja 60 unidades"
$ws.Range("O78").Value = "Región de Arica y Parinacota"
$ws.Range("P78").Value = 133
$ws.Range("Q78").Value = 60
$ws.Range("R78").Value = "Hortaliza"
